# fix: add fields for entities (avatarUrl for user, name for page)
#
# The "Sheet1" worksheet (the active sheet, tabSelected) holds a table of
# DB entities and their columns, one entity per row. This change adds two
# missing fields that were left out of the original schema dump:
#   - Users   (row 3)  gains a new "avatarUrl" column at the end of its
#                       field list (K3).
#   - Page    (row 13) gains a new "name" column inserted right after the
#                       primary key, pushing the existing fields
#                       (category_id, owned_by, avatar_url, background_url,
#                       address, mobile, email, instagram) one cell to the
#                       right (D13:K13 -> E13:L13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Users: add "avatarUrl" as a new trailing field (K3) -------------------
# K3 was a blank placeholder cell (style s=5, matching the other untouched
# placeholder columns); once it holds real field data it should pick up the
# same "filled field" formatting (style s=4) used by its neighbour J3.
$ws.Range("J3").Copy()
$ws.Range("K3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("K3").Value = "avatarUrl"

# --- Page: insert "name" field, shifting later fields one column right ----
# L13 is currently the empty placeholder right after the last used field
# (K13); it needs the "filled field" style before it receives a value, so
# copy formatting from K13 first.
$ws.Range("K13").Copy()
$ws.Range("L13").PasteSpecial(-4122)  # xlPasteFormats

# Shift D13:K13 -> E13:L13 (right to left, so nothing is clobbered).
$ws.Range("L13").Value = $ws.Range("K13").Value2
$ws.Range("K13").Value = $ws.Range("J13").Value2
$ws.Range("J13").Value = $ws.Range("I13").Value2
$ws.Range("I13").Value = $ws.Range("H13").Value2
$ws.Range("H13").Value = $ws.Range("G13").Value2
$ws.Range("G13").Value = $ws.Range("F13").Value2
$ws.Range("F13").Value = $ws.Range("E13").Value2
$ws.Range("E13").Value = $ws.Range("D13").Value2

# New field inserted at the front of the Page column list.
$ws.Range("D13").Value = "name"

# --- Match the author's final selection on save -----------------------
[void]$ws.Range("K3").Select()
